$wb = $excel.ActiveWorkbook
$new = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$new.Name = "Feuil4"

$rng = $new.Range("B6:D11")
$rng.Borders.LineStyle = 1
$rng.Borders.Color = 0
Write-Host "done1"
